# edit.ps1 - applies the commit's changes to the resume document.
$d = $word.ActiveDocument

function Replace-Text($doc, $old, $new) {
    $doc.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

function Get-ParaIndexContaining($doc, $text) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($text)) {
            return $i
        }
    }
    return -1
}

function Remove-Paragraphs($doc, $startText, $endText) {
    $startIdx = Get-ParaIndexContaining $doc $startText
    $endIdx = Get-ParaIndexContaining $doc $endText
    $p1 = $doc.Paragraphs.Item($startIdx)
    $p2 = $doc.Paragraphs.Item($endIdx)
    $r = $doc.Range($p1.Range.Start, $p2.Range.End)
    $r.Delete()
}

# 1. Header name
Replace-Text $d "DHEERAJ CHAND" "Dheeraj Chand"

# 2. Subtitle / professional title
Replace-Text $d "Senior Geospatial Data Engineer & Technical Architect" "Professional Title"

# 3. Contact info (phone + email formatting)
Replace-Text $d "(202) 550-7110 | Dheeraj.Chand@gmail.com" "202.550.7110 | dheeraj.chand@gmail.com"

# 4. Professional summary: years of experience + drop product names
Replace-Text $d "Senior Data Engineer with 20+ years of expertise" "Senior Data Engineer with 21 years of expertise"
Replace-Text $d "Proven track record architecting production systems like BALLISTA and DAMON serving thousands of users" "Proven track record architecting production systems serving thousands of users"

# 5. Company name anonymized
Replace-Text $d "Siege Analytics, Austin, TX | 2005 " "Your Company Name, Your City, ST | 2005 "

# 6. Bullet: drop named platforms
Replace-Text $d "▸ Architected and engineered production geospatial platforms including BALLISTA (redistricting) and DAMON (boundary estimation) serving thousands of analysts" "▸ Architected and engineered production geospatial platforms serving thousands of analysts"

# 7. Bullet: drop "campaign finance"
Replace-Text $d "▸ Implemented fraud detection systems processing multi-terabyte campaign finance datasets with real-time spatial analysis capabilities" "▸ Implemented fraud detection systems processing multi-terabyte datasets with real-time spatial analysis capabilities"

# 8. Remove the DATA PRODUCTS MANAGER / ANALYTICS SUPERVISOR / SOFTWARE ENGINEER / SENIOR DATA ANALYST
#    job blocks entirely (everything between the PARTNER block and KEY ACHIEVEMENTS heading)
Remove-Paragraphs $d "DATA PRODUCTS MANAGER" "Led technical evaluation of 1,200+ vendor proposals"

# 9. Achievements bullets: drop product names
Replace-Text $d "✓ Architected BALLISTA redistricting platform processing Census data for thousands of analysts with real-time PostGIS collaborative editing" "✓ Architected redistricting platform processing Census data for thousands of analysts with real-time PostGIS collaborative editing"
Replace-Text $d "✓ Built DAMON boundary estimation system using advanced PostGIS algorithms and incomplete data without machine learning requirements" "✓ Built boundary estimation system using advanced PostGIS algorithms and incomplete data without machine learning requirements"

# 10. Replace SimCrisis bullet text (this paragraph survives; subsequent ones are removed)
Replace-Text $d "✓ Developed SimCrisis geospatial simulation platform integrating NetLogo multi-agent modeling with GeoDjango web interface" "✓ Developed geospatial simulation platform integrating multi-agent modeling with web interface"

# 11. Remove remaining achievement bullets/sections (RACSO bullet through end of document content)
Remove-Paragraphs $d "Created production-scale survey platform RACSO" "Developed comprehensive spatial data governance frameworks ensuring quality across petabyte-scale geospatial warehouses"

Write-Output "edit complete"
